# Generate Report for Handback
# Update the timestamps recorded for the 9ad94a91-b3c8-4655-835a-0b7ba4c8eeda
# handback entry across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for the 9ad94a91... row
$wsOverview.Range("G3").Value = "2016-08-15 12:42:17"

# zh-cn: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsZhCn.Range("H3").Value = "2016-08-15 12:42:13"
$wsZhCn.Range("K3").Value = "2016-08-15 12:42:31"

# de-de: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsDeDe.Range("H3").Value = "2016-08-15 12:42:17"
$wsDeDe.Range("K3").Value = "2016-08-15 12:42:37"
